$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the changed data values on row 3
$ws.Range("E3").Value = 4
$ws.Range("G3").Value = -3
$ws.Range("H3").Value = 13

# Move the active selection from I1 to E3, matching the saved sheet view
$ws.Range("E3").Select()
